# aggiornamento fino a 20/09/2021
# Appends rows 375-385 (dates 2021-09-10 .. 2021-09-20) to the COVID-style
# "nuovi pos." report on Sheet1, extending the data range from A1:D374 to
# A1:D385.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date serial, nuovi pos. (B), somma mobile 7gg. (C), somma mobile 7gg. per 100mila ab. (D)
$data = @(
    @(44449, 0, 9, 52.52407353370295),
    @(44450, 1, 5, 29.18004085205719),
    @(44451, 1, 6, 35.01604902246864),
    @(44452, 0, 3, 17.50802451123432),
    @(44453, 2, 5, 29.18004085205719),
    @(44454, 0, 5, 29.18004085205719),
    @(44455, 2, 6, 35.01604902246864),
    @(44456, 2, 8, 46.68806536329151),
    @(44457, 0, 7, 40.85205719288007),
    @(44458, 2, 8, 46.68806536329151),
    @(44459, 1, 9, 52.52407353370295)
)

$startRow = 375
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

$endRow = $startRow + $data.Count - 1

# Column A carries the "date" style (s="2", numFmt YYYY-MM-DD HH:MM:SS) used
# by every other row in the column - copy it down from the last existing row
# so the new cells pick up the same cellXf instead of minting a new one.
$ws.Range("A374").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
